$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text format first so numeric-looking
# strings like "418.30" or "10.00" are stored as literal text (matching
# the source data's inlineStr type) instead of being coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.768.54"
$ws.Range("D3").Value = "3.601.14"
$ws.Range("D5").Value = "418.30"
$ws.Range("D6").Value = "130.84"
$ws.Range("D8").Value = "3.592.26"
$ws.Range("D10").Value = "0.769"
$ws.Range("D13").Value = "42.84"
$ws.Range("D14").Value = "10.00"
$ws.Range("D15").Value = "4.151.44"
$ws.Range("D18").Value = "3.598.08"
$ws.Range("D19").Value = "1.15"
$ws.Range("D20").Value = "67.715.15"
$ws.Range("D21").Value = "12.38"
$ws.Range("D22").Value = "464.97"
$ws.Range("D23").Value = "88.70"
$ws.Range("D24").Value = "3.14"
$ws.Range("D25").Value = "13.41"
$ws.Range("D26").Value = "10.31"
$ws.Range("D27").Value = "3.38"
$ws.Range("D28").Value = "36.48"
$ws.Range("D29").Value = "4.87"
$ws.Range("D31").Value = "12.46"
$ws.Range("D32").Value = "7.48"
$ws.Range("D35").Value = "41.48"
$ws.Range("D37").Value = "56.88"
$ws.Range("D39").Value = "0.0₃0720"
$ws.Range("D43").Value = "148.76"
$ws.Range("D44").Value = "2.74"
$ws.Range("D46").Value = "4.34"
$ws.Range("D47").Value = "0.311"
$ws.Range("D49").Value = "2.35"
$ws.Range("D51").Value = "15.74"

# Revert to the default (unstyled) cell style now that the values are
# committed as text, so no extra style index lingers on these cells.
$ws.Range("D2:D51").Style = "Normal"

$ws.Range("E2").Value = "  +7.93%  "
$ws.Range("E3").Value = "  +3.80%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  +3.61%  "
$ws.Range("E8").Value = "  +3.77%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("E10").Value = "  +5.17%  "
$ws.Range("E11").Value = "  +18.28%  "
$ws.Range("E12").Value = "  +54.62%  "
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("E14").Value = "  +2.48%  "
$ws.Range("E15").Value = "  +2.93%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("E18").Value = "  +3.47%  "
$ws.Range("E19").Value = "  +5.40%  "
$ws.Range("E20").Value = "  +7.73%  "
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("E23").Value = "  -2.28%  "
$ws.Range("E24").Value = "  -4.84%  "
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("E28").Value = "  +8.91%  "
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("E30").Value = "  +3.97%  "
$ws.Range("E31").Value = "  +2.40%  "
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("E33").Value = "  +4.36%  "
$ws.Range("E34").Value = "  -3.16%  "
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -2.42%  "
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("E39").Value = "  +25.69%  "
$ws.Range("E40").Value = "  +8.01%  "
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("E45").Value = "  -2.06%  "
$ws.Range("E46").Value = "  -3.23%  "
$ws.Range("E47").Value = "  -3.38%  "
$ws.Range("E48").Value = "  -3.79%  "
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("E50").Value = "  +16.80%  "
$ws.Range("E51").Value = "  -4.09%  "
